# Weekly update: a new daily price record for "Ají" (Inferno / Primera,
# Región de Arica y Parinacota, $/caja 12 kilos) is inserted at row 259,
# pushing the existing rows 259-322 down to 260-323.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 259, shifting rows 259:322 down
# to 260:323 (dimension grows from A1:R322 to A1:R323).
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new observation.
$ws.Range("A259").Value = 9
$ws.Range("B259").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C259").Value = "Metropolitana"
$ws.Range("D259").Value = 44782
$ws.Range("E259").Value = 13
$ws.Range("F259").Value = 100112021
$ws.Range("G259").Value = "Ají"
$ws.Range("H259").Value = "Inferno"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 25
$ws.Range("K259").Value = 12000
$ws.Range("L259").Value = 13000
$ws.Range("M259").Value = 12480
$ws.Range("N259").Value = "$/caja 12 kilos"
$ws.Range("O259").Value = "Región de Arica y Parinacota"
$ws.Range("P259").Value = 1040
$ws.Range("Q259").Value = 12
$ws.Range("R259").Value = "Hortaliza"
